$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at position 54 ---
$ws.Rows(54).Insert()

$ws.Range("A54").Value = 5
$ws.Range("B54").Value = "Macroferia Regional de Talca"
$ws.Range("C54").Value = "Maule"
$ws.Range("D54").Value = 44720
$ws.Range("D54").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E54").Value = 7
$ws.Range("F54").Value = 100112031
$ws.Range("G54").Value = "Poroto verde"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 150
$ws.Range("K54").Value = 23000
$ws.Range("L54").Value = 23000
$ws.Range("M54").Value = 23000
$ws.Range("N54").Value = "$/saco 25 kilos"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 920
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"

# --- Insert second new row at position 68 ---
# (after the first insertion, the original row 67 now sits at row 68;
#  inserting here pushes it - and everything below - down by one more)
$ws.Rows(68).Insert()

$ws.Range("A68").Value = 5
$ws.Range("B68").Value = "Macroferia Regional de Talca"
$ws.Range("C68").Value = "Maule"
$ws.Range("D68").Value = 44721
$ws.Range("D68").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E68").Value = 7
$ws.Range("F68").Value = 100112031
$ws.Range("G68").Value = "Poroto verde"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 150
$ws.Range("K68").Value = 23000
$ws.Range("L68").Value = 23000
$ws.Range("M68").Value = 23000
$ws.Range("N68").Value = "$/malla 25 kilos"
$ws.Range("O68").Value = "Región de Arica y Parinacota"
$ws.Range("P68").Value = 920
$ws.Range("Q68").Value = 25
$ws.Range("R68").Value = "Hortaliza"
